$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "DO"
$ws.Range("I2").Value = "M3"
$ws.Range("J2").Value = "DO"
$ws.Range("P2").Value = "DO"
$ws.Range("S2").Value = "M3"
$ws.Range("T2").Value = "M1"
$ws.Range("W2").Value = "DO"
$ws.Range("Y2").Value = "M1"
$ws.Range("AC2").Value = "M3"
$ws.Range("B3").Value = "DO"
$ws.Range("D3").Value = "M2"
$ws.Range("H3").Value = "M2"
$ws.Range("I3").Value = "DO"
$ws.Range("K3").Value = "M1"
$ws.Range("L3").Value = "M3"
$ws.Range("N3").Value = "A1"
$ws.Range("O3").Value = "M2"
$ws.Range("P3").Value = "DO"
$ws.Range("Q3").Value = "M1"
$ws.Range("S3").Value = "M1"
$ws.Range("U3").Value = "A1"
$ws.Range("V3").Value = "M3"
$ws.Range("W3").Value = "DO"
$ws.Range("X3").Value = "M3"
$ws.Range("Y3").Value = "M1"
$ws.Range("Z3").Value = "M2"
$ws.Range("AA3").Value = "M2"
$ws.Range("AB3").Value = "A2"
$ws.Range("B4").Value = "DO"
$ws.Range("C4").Value = "M2"
$ws.Range("F4").Value = "M3"
$ws.Range("G4").Value = "M3"
$ws.Range("I4").Value = "DO"
$ws.Range("J4").Value = "M3"
$ws.Range("K4").Value = "M2"
$ws.Range("M4").Value = "M2"
$ws.Range("O4").Value = "M2"
$ws.Range("P4").Value = "DO"
$ws.Range("Y4").Value = "DO"
$ws.Range("AA4").Value = "M1"
$ws.Range("AC4").Value = "M3"
$ws.Range("B5").Value = "DO"
$ws.Range("C5").Value = "M2"
$ws.Range("D5").Value = "M1"
$ws.Range("E5").Value = "M1"
$ws.Range("F5").Value = "M2"
$ws.Range("H5").Value = "M3"
$ws.Range("I5").Value = "DO"
$ws.Range("L5").Value = "M2"
$ws.Range("M5").Value = "M3"
$ws.Range("N5").Value = "A1"
$ws.Range("O5").Value = "M1"
$ws.Range("P5").Value = "DO"
$ws.Range("T5").Value = "M1"
$ws.Range("U5").Value = "M3"
$ws.Range("V5").Value = "M2"
$ws.Range("W5").Value = "DO"
$ws.Range("X5").Value = "M3"
$ws.Range("Y5").Value = "M1"
$ws.Range("Z5").Value = "M2"
$ws.Range("AB5").Value = "A2"
$ws.Range("AC5").Value = "M1"
$ws.Range("B6").Value = "DO"
$ws.Range("E6").Value = "M1"
$ws.Range("I6").Value = "DO"
$ws.Range("K6").Value = "M1"
$ws.Range("L6").Value = "M3"
$ws.Range("N6").Value = "A1"
$ws.Range("O6").Value = "M2"
$ws.Range("P6").Value = "DO"
$ws.Range("Q6").Value = "M1"
$ws.Range("R6").Value = "M3"
$ws.Range("S6").Value = "M2"
$ws.Range("T6").Value = "M1"
$ws.Range("U6").Value = "A2"
$ws.Range("V6").Value = "M2"
$ws.Range("W6").Value = "DO"
$ws.Range("X6").Value = "M1"
$ws.Range("Z6").Value = "M2"
$ws.Range("AA6").Value = "M3"
$ws.Range("AC6").Value = "M3"
$ws.Range("E7").Value = "M3"
$ws.Range("H7").Value = "DO"
$ws.Range("I7").Value = "DO"
$ws.Range("M7").Value = "M3"
$ws.Range("P7").Value = "DO"
$ws.Range("Q7").Value = "M3"
$ws.Range("S7").Value = "M1"
$ws.Range("W7").Value = "DO"
$ws.Range("X7").Value = "M3"
$ws.Range("AA7").Value = "M1"
$ws.Range("B8").Value = "DO"
$ws.Range("D8").Value = "M2"
$ws.Range("F8").Value = "M1"
$ws.Range("H8").Value = "M3"
$ws.Range("I8").Value = "DO"
$ws.Range("M8").Value = "M3"
$ws.Range("N8").Value = "A1"
$ws.Range("O8").Value = "M1"
$ws.Range("P8").Value = "DO"
$ws.Range("Q8").Value = "M1"
$ws.Range("S8").Value = "M2"
$ws.Range("T8").Value = "M2"
$ws.Range("U8").Value = "A2"
$ws.Range("W8").Value = "DO"
$ws.Range("X8").Value = "M2"
$ws.Range("Z8").Value = "M1"
$ws.Range("AA8").Value = "M1"
$ws.Range("AB8").Value = "A2"
$ws.Range("AC8").Value = "M3"
$ws.Range("B9").Value = "DO"
$ws.Range("C9").Value = "M1"
$ws.Range("F9").Value = "M3"
$ws.Range("H9").Value = "M1"
$ws.Range("I9").Value = "DO"
$ws.Range("J9").Value = "M1"
$ws.Range("N9").Value = "A1"
$ws.Range("O9").Value = "M1"
$ws.Range("P9").Value = "DO"
$ws.Range("Q9").Value = "M3"
$ws.Range("R9").Value = "M1"
$ws.Range("S9").Value = "M2"
$ws.Range("T9").Value = "M2"
$ws.Range("V9").Value = "M1"
$ws.Range("W9").Value = "DO"
$ws.Range("X9").Value = "M1"
$ws.Range("Y9").Value = "M3"
$ws.Range("AB9").Value = "A2"
$ws.Range("AC9").Value = "M2"
$ws.Range("B10").Value = "DO"
$ws.Range("D10").Value = "M1"
$ws.Range("F10").Value = "M3"
$ws.Range("G10").Value = "A1"
$ws.Range("H10").Value = "M1"
$ws.Range("I10").Value = "DO"
$ws.Range("L10").Value = "M2"
$ws.Range("M10").Value = "M2"
$ws.Range("N10").Value = "A2"
$ws.Range("O10").Value = "M3"
$ws.Range("P10").Value = "DO"
$ws.Range("R10").Value = "M2"
$ws.Range("T10").Value = "M2"
$ws.Range("U10").Value = "A2"
$ws.Range("V10").Value = "M3"
$ws.Range("Z10").Value = "M3"
$ws.Range("AB10").Value = "DO"
$ws.Range("AC10").Value = "M1"
